$wb = $excel.ActiveWorkbook

# Rename sheets
$wsConn = $wb.Worksheets.Item("CONEXIONES CON")
$wsConn.Name = "CONNECTIVITY"

$wsCoord = $wb.Worksheets.Item("COORDENADAS COOR")
$wsCoord.Name = "COORDINATES"

$wsFree = $wb.Worksheets.Item("NODOS LIBRES NL")
$wsFree.Name = "FREE NODES"

# Update formulas on FREE NODES sheet (now referencing COORDINATES unquoted, since no spaces)
for ($r = 1; $r -le 12; $r++) {
    $wsFree.Range("A$r").Formula = "=+COORDINATES!A$r"
}

# Selection / active sheet changes
$wsCoord.Range("B7").Select()
$wsConn.Activate()
$wsConn.Range("C6").Select()
